# Update the crypto price list with the latest scraped values.
# GitHub Actions symbol-list refresh.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Column D ("Price") values are numeric-looking text; keep them stored
#     as text (matching the original inlineStr type) by forcing the cell
#     number format to Text before writing the new value. ---
$priceUpdates = @{
    2  = "243.91"
    3  = "23.20"
    4  = "5.420"
    5  = "0.05973"
    6  = "3.431"
    7  = "6.526"
    8  = "0.8110"
    9  = "0.9311"
    10 = "0.1431"
    11 = "0.07386"
    12 = "0.03300"
    13 = "0.03076"
    14 = "0.09347"
    15 = "3.869"
    16 = "0.001580"
    17 = "0.04700"
    18 = "0.0005891"
    19 = "0.005969"
    20 = "0.001271"
    21 = "0.004916"
    22 = "0.00006801"
    23 = "3.575"
    27 = "0.0003700"
    40 = "0.03959"
    41 = "0.1079"
    42 = "0.002571"
    43 = "0.003077"
    44 = "0.009424"
    45 = "0.00005221"
    46 = "0.00000000750"
    47 = "0.7251"
    48 = "0.002306"
    49 = "0.00002100"
    50 = "0.0002000"
}

foreach ($row in $priceUpdates.Keys) {
    $cell = $ws.Cells.Item($row, 4)
    $cell.NumberFormat = "@"
    $cell.Value = $priceUpdates[$row]
}

# --- Rows 41-43: the three coins were reshuffled (Kick/BKEX/CEJI rotated)
#     and their linked URL + rank label columns need to follow suit. ---
$ws.Cells.Item(41, 2).Value = "BKEXToken"
$ws.Cells.Item(41, 3).Value = "https://coinranking.com/coin/IPeThtYgk+bkextoken-bkk"
$ws.Cells.Item(41, 5).Value = "40BKEXTokenBKK"

$ws.Cells.Item(42, 2).Value = "CEJI"
$ws.Cells.Item(42, 3).Value = "https://coinranking.com/coin/SbKjCVJCh+ceji-ceji"
$ws.Cells.Item(42, 5).Value = "41CEJICEJI"

$ws.Cells.Item(43, 2).Value = "KickToken"
$ws.Cells.Item(43, 3).Value = "https://coinranking.com/coin/F_Yv9Cu7pPL3Y+kicktoken-kick"
$ws.Cells.Item(43, 5).Value = "42KickTokenKICKWorstin24h"

# --- Row 47: CoinbaseStockToken is no longer flagged as the 24h worst. ---
$ws.Cells.Item(47, 5).Value = "46CoinbaseStockTokenCOIN"
